$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A42").Value = "2025/12/04 11:00"
$ws.Range("B42").Value = "-"
$ws.Range("C42").Value = "-"
$ws.Range("D42").Value = "-"
$ws.Range("E42").Value = "-"
$ws.Range("F42").Value = "-"
$ws.Range("G42").Value = "-"
